# Actualización automática de scrims_actualizado.xlsx (2025-07-24 02:20:00)
# Adds newly-submitted scrim rows to three sheets:
#   - "Hideout"        : append row 5  (duplicate of row 4, new timestamp)
#   - "Crystal Arcade"  : append rows 34-35 (new submissions)
#   - "Hard Rock Mine"  : append rows 13-14 (new submissions)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Hideout": A3:N4 -> A3:N5
# New row 5 duplicates row 4's player roster/team, only the
# Timestamp (column N) differs.
# ---------------------------------------------------------------
$wsHideout = $wb.Worksheets.Item("Hideout")
$wsHideout.Range("A4:N4").Copy($wsHideout.Range("A5:N5"))
$wsHideout.Range("N5").Value = "20250723T234900.000Z"

# ---------------------------------------------------------------
# Sheet "Crystal Arcade": A3:N33 -> A3:N35
# Two new submissions (rows 34 & 35) for the same matchup.
# ---------------------------------------------------------------
$wsCrystal = $wb.Worksheets.Item("Crystal Arcade")

# Row 34
$wsCrystal.Range("A4:N4").Copy($wsCrystal.Range("A34:N34"))
$wsCrystal.Range("A34").Value = "TARA"
$wsCrystal.Range("B34").Value = "CORDELIUS"
$wsCrystal.Range("C34").Value = "GUS"
$wsCrystal.Range("D34").Value = "NITA"
$wsCrystal.Range("E34").Value = "DRACO"
$wsCrystal.Range("F34").Value = "FINX"
$wsCrystal.Range("G34").Value = "Equipo 1"
$wsCrystal.Range("H34").Value = "Always cool🧃"
$wsCrystal.Range("I34").Value = "Dan-YT"
$wsCrystal.Range("J34").Value = "DMO|Marco"
$wsCrystal.Range("K34").Value = "SKC|Rhz"
$wsCrystal.Range("L34").Value = "Jxcccr🐻‍❄️"
$wsCrystal.Range("M34").Value = "SKC|Kr"
$wsCrystal.Range("N34").Value = "20250723T235810.000Z"

# Row 35 (same roster/team as row 34, different timestamp)
$wsCrystal.Range("A34:N34").Copy($wsCrystal.Range("A35:N35"))
$wsCrystal.Range("N35").Value = "20250723T235543.000Z"

# ---------------------------------------------------------------
# Sheet "Hard Rock Mine": A3:N12 -> A3:N14
# Two new submissions (rows 13 & 14) for the same matchup.
# ---------------------------------------------------------------
$wsHardRock = $wb.Worksheets.Item("Hard Rock Mine")

# Row 13
$wsHardRock.Range("A6:N6").Copy($wsHardRock.Range("A13:N13"))
$wsHardRock.Range("A13").Value = "JAE-YONG"
$wsHardRock.Range("B13").Value = "GUS"
$wsHardRock.Range("C13").Value = "HANK"
$wsHardRock.Range("D13").Value = "JANET"
$wsHardRock.Range("E13").Value = "BO"
$wsHardRock.Range("F13").Value = "SHADE"
$wsHardRock.Range("G13").Value = "Equipo 1"
$wsHardRock.Range("H13").Value = "Solar Ray ☀️"
$wsHardRock.Range("I13").Value = "Finki is back."
$wsHardRock.Range("J13").Value = "Xyz"
$wsHardRock.Range("K13").Value = "BC*|Jubileubr"
$wsHardRock.Range("L13").Value = "LOUD|Edinho"
$wsHardRock.Range("M13").Value = "CASA|Doritos"
$wsHardRock.Range("N13").Value = "20250724T001158.000Z"

# Row 14 (same roster/team as row 13, different timestamp)
$wsHardRock.Range("A13:N13").Copy($wsHardRock.Range("A14:N14"))
$wsHardRock.Range("N14").Value = "20250724T001013.000Z"
